$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text with the new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$rangeA1 = $wsHoja1.Range("A1")
$oldText = $rangeA1.Text
$newText = $oldText.Replace("1000 Bs = 3.25 = 12441.6 pesos", "1000 Bs = 3.23 = 12334.92 pesos")
$newText = $newText.Replace("12441.6 pesos = 3.23 = 968.5 Bs", "12334.92 pesos = 3.21 = 967.77 Bs")
$rangeA1.Value = $newText

# --- Update tasas sheet rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 309.69
$wsTasas.Range("O10").Value = 3820
$wsTasas.Range("N12").Value = 3839
$wsTasas.Range("O12").Value = 301.2
